# Updated cryptos list on Sun Sep 17 12:44:03 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "26.791.86"; E = "  +0.47%  " },
    @{ Row = 3;  D = "1.644.93";  E = "  +0.08%  " },
    @{ Row = 4;  E = "  +0.71%  " },
    @{ Row = 5;  D = "216.79";    E = "  +0.72%  " },
    @{ Row = 6;  E = "  -0.47%  " },
    @{ Row = 7;  E = "  +0.33%  " },
    @{ Row = 8;  E = "  -0.35%  " },
    @{ Row = 9;  D = "0.0629";    E = "  +0.30%  " },
    @{ Row = 10; D = "19.22";     E = "  -0.10%  " },
    @{ Row = 11; D = "0.0845";    E = "  +0.46%  " },
    @{ Row = 12; D = "1.641.78";  E = "  -0.47%  " },
    @{ Row = 13; E = "  -0.70%  " },
    @{ Row = 14; E = "  -0.19%  " },
    @{ Row = 15; D = "64.74";     E = "  -0.59%  " },
    @{ Row = 16; D = "26.793.00"; E = "  +0.38%  " },
    @{ Row = 17; E = "  -1.25%  " },
    @{ Row = 18; D = "214.75";    E = "  -0.73%  " },
    @{ Row = 19; E = "  +0.51%  " },
    @{ Row = 20; D = "4.39";      E = "  +0.97%  " },
    @{ Row = 21; D = "2.45";      E = "  +8.26%  " },
    @{ Row = 22; E = "  -0.03%  " },
    @{ Row = 23; D = "9.33";      E = "  -1.90%  " },
    @{ Row = 24; D = "146.02";    E = "  +0.06%  " },
    @{ Row = 25; E = "  +0.52%  " },
    @{ Row = 26; E = "  -1.15%  " },
    @{ Row = 27; E = "  +0.41%  " },
    @{ Row = 28; D = "15.67";     E = "  -0.28%  " },
    @{ Row = 29; D = "0.0510";    E = "  -1.39%  " },
    @{ Row = 30; E = "  +0.71%  " },
    @{ Row = 31; D = "3.36";      E = "  -0.73%  " },
    @{ Row = 32; D = "3.00";      E = "  -1.42%  " },
    @{ Row = 33; D = "1.285.91";  E = "  +0.78%  " },
    @{ Row = 34; E = "  -0.30%  " },
    @{ Row = 35; D = "2.44";      E = "  +1.37%  " },
    @{ Row = 36; E = "  -1.32%  " },
    @{ Row = 37; D = "0.536";     E = "  +0.23%  " },
    @{ Row = 38; E = "  -1.26%  " },
    @{ Row = 40; E = "  -1.10%  " },
    @{ Row = 41; D = "2.24";      E = "  -0.32%  " },
    @{ Row = 42; D = "5.33";      E = "  -2.35%  " },
    @{ Row = 43; D = "1.785.29";  E = "  +0.10%  " },
    @{ Row = 44; D = "61.73";     E = "  +3.62%  " },
    @{ Row = 45; D = "91.90";     E = "  +0.60%  " },
    @{ Row = 46; D = "1.60";      E = "  -0.13%  " },
    @{ Row = 47; D = "0.0518";    E = "  +0.40%  " },
    @{ Row = 48; D = "7.65";      E = "  -1.94%  " },
    @{ Row = 49; D = "0.0970";    E = "  +0.20%  " },
    @{ Row = 50; D = "0.408";     E = "  +0.42%  " },
    @{ Row = 51; E = "  +0.39%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        # Column D holds textual "price" figures (e.g. "26.791.86",
        # "216.79"). Several of these look like plain numbers to Excel's
        # auto-detection, which would silently convert them into numeric
        # cells (losing the original text form / trailing zeros). Force
        # the cell to Text before writing, then restore the default
        # "Normal" style so no stray formatting is left behind.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        # Column E values are padded with spaces ("  +0.47%  "), which
        # keeps Excel from reinterpreting them as numeric/percentage
        # values, so a plain assignment is safe here.
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
